# Update column G ("K") values for rows 2-7 on Sheet1
# These values represent strikeouts recalculated using K instead of Strike#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 7
$ws.Range("G3").Value = 5
$ws.Range("G4").Value = 10
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
